$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in columns D and E store plain text (e.g. "46.578.48", "  +1.22%  ").
# Some D values (like "305.46") look numeric to Excel, so force a text
# NumberFormat before assignment, then restore the default "Normal" style
# so the cell style matches the original (unstyled) cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.578.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.578.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +9.39%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.09%  "

$ws.Range("E7").Value = "  +5.37%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.573"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0838"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.87%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.973.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.32%  "

$ws.Range("E14").Value = "  +1.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.585.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.905"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +11.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +9.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "46.721.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.63%  "

$ws.Range("E20").Value = "  +3.90%  "

$ws.Range("E21").Value = "  +10.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.21%  "

$ws.Range("E24").Value = "  +3.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +35.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.98%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.91%  "

$ws.Range("E33").Value = "  +23.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.72%  "

$ws.Range("E35").Value = "  +6.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "150.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.38%  "

$ws.Range("E37").Value = "  +3.77%  "

$ws.Range("E38").Value = "  +4.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.85%  "

$ws.Range("E41").Value = "  +12.02%  "

$ws.Range("E42").Value = "  +7.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.019.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +26.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("E47").Value = "  +0.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "108.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.201"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.833.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.10%  "
